$d = $word.ActiveDocument

$pairs = @(
    @("39×96=3744", "18×55=990"),
    @("23×21=483", "58×63=3654"),
    @("38×97=3686", "35×77=2695"),
    @("37×94=3478", "39×83=3237"),
    @("77×82=6314", "28×47=1316"),
    @("20×57=1140", "33×77=2541"),
    @("49×26=1274", "42×64=2688"),
    @("32×31=992", "13×87=1131"),
    @("63×93=5859", "29×79=2291"),
    @("47×45=2115", "31×79=2449"),
    @("34×37=1258", "14×49=686"),
    @("28×98=2744", "20×87=1740"),
    @("31×21=651", "24×96=2304"),
    @("97×96=9312", "53×26=1378"),
    @("48×42=2016", "33×56=1848"),
    @("78×80=6240", "49×95=4655"),
    @("89×50=4450", "56×96=5376"),
    @("88×52=4576", "15×80=1200"),
    @("56×61=3416", "76×19=1444"),
    @("40×20=800", "40×56=2240"),
    @("19×78=1482", "28×76=2128"),
    @("26×62=1612", "34×84=2856"),
    @("66×35=2310", "60×78=4680"),
    @("91×30=2730", "81×89=7209"),
    @("99×45=4455", "48×55=2640")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
